$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Range("B2").Value = 0.3354314908498832
$ws1.Range("B3").Value = 59.83554148444802
$ws1.Range("B4").Value = 61.58417617154723
$ws1.Range("B5").Value = 62.22895309416919
$ws1.Range("B6").Value = 62.89503060659936
$ws1.Range("B7").Value = 62.86109488146062
$ws1.Range("B8").Value = 62.76840445400494
$ws1.Range("B9").Value = 62.81774446234635
$ws1.Range("B10").Value = 63.77283542464308
$ws1.Range("B11").Value = 63.94162805970876
$ws1.Range("B12").Value = 64.6945156943917
$ws1.Range("B13").Value = 64.71609285619667
$ws1.Range("B14").Value = 64.51049183329771
$ws1.Range("B15").Value = 64.31153458529798
$ws1.Range("B16").Value = 64.64129683832219
$ws1.Range("B17").Value = 64.61748084253712
$ws1.Range("B18").Value = 64.69786712441247
$ws1.Range("B19").Value = 64.04407180430012
$ws1.Range("B20").Value = 64.6534359094167
$ws1.Range("B21").Value = 64.79706615706176
$ws1.Range("B22").Value = 64.40488823231416
$ws1.Range("B23").Value = 66.02539923867812
$ws1.Range("B24").Value = 64.96151883986936
$ws1.Range("B25").Value = 65.49741265638885
$ws1.Range("B26").Value = 66.58590172946425
$ws1.Range("B27").Value = 66.02171152035156
$ws1.Range("B28").Value = 67.56585866372433

$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Range("B2").Value = 0.2004418719579071
$ws3.Range("B3").Value = 59.70055186555603
$ws3.Range("B4").Value = 61.44918655265525
$ws3.Range("B5").Value = 62.09396347527721
$ws3.Range("B6").Value = 62.76004098770738
$ws3.Range("B7").Value = 62.72610526256864
$ws3.Range("B8").Value = 62.63341483511296
$ws3.Range("B9").Value = 62.68275484345437
$ws3.Range("B10").Value = 63.6378458057511
$ws3.Range("B11").Value = 63.80663844081678
$ws3.Range("B12").Value = 64.55952607549972
$ws3.Range("B13").Value = 64.58110323730469
$ws3.Range("B14").Value = 64.37550221440573
$ws3.Range("B15").Value = 64.17654496640601
$ws3.Range("B16").Value = 64.50630721943021
$ws3.Range("B17").Value = 64.48249122364514
$ws3.Range("B18").Value = 64.56287750552049
$ws3.Range("B19").Value = 63.90908218540815
$ws3.Range("B20").Value = 64.51844629052472
$ws3.Range("B21").Value = 64.66207653816979
$ws3.Range("B22").Value = 64.26989861342219
$ws3.Range("B23").Value = 65.89040961978614
$ws3.Range("B24").Value = 64.82652922097738
$ws3.Range("B25").Value = 65.36242303749687
$ws3.Range("B26").Value = 66.45091211057228
$ws3.Range("B27").Value = 65.88672190145958
$ws3.Range("B28").Value = 67.43086904483235
$ws3.Range("B29").Value = 68.42177158269732
$ws3.Range("B30").Value = 67.21497139880758
$ws3.Range("B31").Value = 66.85768522037939
$ws3.Range("B32").Value = 66.26515782507337
$ws3.Range("B33").Value = 66.48125036682833
